$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the title
# (Heading1) paragraph. We clone paragraph 50 - which already has the exact
# target run layout (empty run, then a bold run) - via copy/paste so the new
# paragraph naturally ends up with the same <w:r/><w:r><w:rPr><w:b/>...
# structure, then we swap in the real text.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$sourcePara = $d.Paragraphs(50)
$sourcePara.Range.Copy()

$titlePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"
$d.Paragraphs(2).Range.Paste()

# Turn the pasted "Play Crystal Forest HD for Free - WMS Online Slot Review"
# (bold) text into "Meta description" (still bold), scoped to this paragraph
# only so the real title in paragraph 1 is left untouched.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Find.Execute("Play Crystal Forest HD for Free - WMS Online Slot Review", $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

# Append the (non-bold) description text after "Meta description".
$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$insertPoint = $d.Range($metaRange.End - 1, $metaRange.End - 1)
$insertPoint.InsertAfter(": Get familiar with the tricks and features of WMS slot game Crystal Forest HD. Play for free and enjoy enchanting graphics and cascading reels.")

# ---------------------------------------------------------------------------
# Change 2: delete the duplicate bold "Play Crystal Forest HD for Free - WMS
# Online Slot Review" paragraph that used to sit right before the final
# italic paragraph (now shifted down one slot to 51 because of Change 1).
# ---------------------------------------------------------------------------
$dupIndex = $d.Paragraphs.Count - 1
$d.Paragraphs($dupIndex).Range.Delete()

# ---------------------------------------------------------------------------
# Change 3: replace the text of the final (italic) paragraph with the new
# DALL-E prompt, keeping the italic run formatting and straight quotes intact
# (avoid Find/Replace's smart-quote substitution by assigning .Text directly
# on a range that excludes the trailing paragraph mark).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "Dear DALLE, I need you to create a feature image for the online slot game ""Crystal Forest HD"". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be portrayed in a whimsical and magical way to reflect the enchanted world of the game. Please ensure that the image is visually engaging and eye-catching to capture the attention of the audience. The color scheme should be in line with the game's graphics, especially with regards to the blue background. I'm confident that you can come up with a fantastic feature image that perfectly aligns with the game's theme and enhances its overall appeal. Thank you in advance for your hard work and creativity. Best regards, [Your Name]"

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
Write-Output ("Para 2: " + $d.Paragraphs(2).Range.Text)
Write-Output ("Last para: " + $d.Paragraphs($d.Paragraphs.Count).Range.Text)
